$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: Shaedon Sharpe, SG,SF, Portland Trail Blazers -> LaMelo Ball, PG,SG, Charlotte Hornets
$ws.Range("A5").Value = "LaMelo Ball"
$ws.Range("B5").Value = "PG,SG"
$ws.Range("C5").Value = "Charlotte Hornets"

# Row 12: Grant Williams, PF,C, Charlotte Hornets -> Bam Adebayo, C, Miami Heat
$ws.Range("A12").Value = "Bam Adebayo"
$ws.Range("B12").Value = "C"
$ws.Range("C12").Value = "Miami Heat"

# Row 13: Bam Adebayo, C, Miami Heat -> Mason Plumlee, C, Phoenix Suns
$ws.Range("A13").Value = "Mason Plumlee"
$ws.Range("B13").Value = "C"
$ws.Range("C13").Value = "Phoenix Suns"

# Row 14: LaMelo Ball, PG,SG, Charlotte Hornets -> Grant Williams, PF,C, Charlotte Hornets
$ws.Range("A14").Value = "Grant Williams"
$ws.Range("B14").Value = "PF,C"
$ws.Range("C14").Value = "Charlotte Hornets"
